$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2:G21").Value = "2021/2022"
$ws.Range("G8").Select() | Out-Null
